$d = $word.ActiveDocument

# --- Step 1: merge the ")" and "." runs at the end of the document into a
# single run, dropping the "_GoBack" bookmark that currently sits between
# them (Word will re-anchor that bookmark at the new edit point instead).
$closeParen = $d.Content
$found = $closeParen.Find.Execute(")", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $mergeRange = $d.Range($closeParen.Start, $closeParen.Start + 2)
    $mergedXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $mergeRange.InsertXML($mergedXml)
}

# --- Step 2: insert a new, empty paragraph at the very start of the
# document holding only the "_GoBack" bookmark (this is where Word actually
# records the last edit position once the change above is made).
$start = $d.Range(0, 0)
$bookmarkXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$start.InsertXML($bookmarkXml) | Out-Null
